$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; existing rows 26-35 shift down to 27-36
$ws.Rows.Item(26).Insert()

# Fill in the new row 26 (Stephan Heunis)
$ws.Range("A26").Value = "Stephan Heunis"
$ws.Range("B26").Value = 25
$ws.Range("B26").Interior.Color = 65535

# Renumber column B for the shifted rows (27-36) to keep the running count
$ws.Range("B27").Value = 26
$ws.Range("B28").Value = 27
$ws.Range("B29").Value = 28
$ws.Range("B30").Value = 29
$ws.Range("B31").Value = 30
$ws.Range("B32").Value = 31
$ws.Range("B33").Value = 32
$ws.Range("B34").Value = 33
$ws.Range("B35").Value = 34
$ws.Range("B36").Value = 35

# Update the selection to match the edited workbook's saved view state
$ws.Range("B36").Select()
